$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as numeric but must stay stored as TEXT
# (matches source formatting of price column, e.g. '0.906').
# Mark them with a Text number format, assign, then drop back to the
# Normal cell style so no '@' format sticks to the cell afterwards.
$textForceCells = @("D5", "D9", "D10", "D11", "D15", "D16", "D18", "D20", "D22", "D23", "D25", "D26", "D28", "D30", "D34", "D37", "D38", "D40", "D42", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.043.49'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '1.660.66'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '216.21'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +2.78%  '
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('D10').Value = '20.15'
$ws.Range('E10').Value = '  +4.81%  '
$ws.Range('D11').Value = '0.0885'
$ws.Range('E11').Value = '  +4.44%  '
$ws.Range('D12').Value = '1.892.96'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').Value = '1.661.71'
$ws.Range('E13').Value = '  +2.78%  '
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').Value = '65.71'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').Value = '27.056.24'
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '236.64'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = '0.0₃0738'
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '7.76'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = '4.44'
$ws.Range('E22').Value = '  +3.46%  '
$ws.Range('D23').Value = '2.25'
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('E24').Value = '  +2.67%  '
$ws.Range('D25').Value = '145.28'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('D26').Value = '7.14'
$ws.Range('E26').Value = '  +2.01%  '
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('D28').Value = '15.91'
$ws.Range('E28').Value = '  +2.54%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '0.0498'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').Value = '1.550.25'
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('D34').Value = '3.09'
$ws.Range('E34').Value = '  +4.24%  '
$ws.Range('E35').Value = '  +6.82%  '
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.906'
$ws.Range('E37').Value = '  +9.27%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.577'
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('E39').Value = '  +2.37%  '
$ws.Range('D40').Value = '6.07'
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '66.70'
$ws.Range('E42').Value = '  +8.61%  '
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('D44').Value = '0.970'
$ws.Range('E44').Value = '  +6.27%  '
$ws.Range('D45').Value = '1.803.36'
$ws.Range('E45').Value = '  +2.62%  '
$ws.Range('D46').Value = '0.778'
$ws.Range('E46').Value = '  +2.16%  '
$ws.Range('D47').Value = '90.67'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('E48').Value = '  +2.54%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.100'
$ws.Range('E49').Value = '  +4.43%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0507'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.70'
$ws.Range('E51').Value = '  +2.82%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
